# Upgrade left table until Javakheti:
#  - rename the sheet to "Terjola"
#  - replace the confidential/unavailable numeric cells in the Urban and
#    Rural rows with placeholder dots ("..." / "…"), matching the Total row
#    which already uses that convention
#  - remove the now-empty row between the data table and the footnote

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Terjola"

# Row 6 ("Urban"): years 2010-2015 (B:G) become "..." placeholders;
# years 2016-2023 (H:O) become the "…" placeholder already used elsewhere.
$ws.Range("B6:G6").Value = "..."
$ws.Range("H6:O6").Value = "…"

# Row 7 ("Rural"): years 2011-2014 (C:F) become "..." placeholders;
# years 2016-2023 (H:O) become "…" placeholders. B7 (2010) and G7 (2015)
# keep their existing numbers.
$ws.Range("C7:F7").Value = "..."
$ws.Range("H7:O7").Value = "…"

# Remove the blank row separating the table from the footnote so the
# footnote moves from row 9 up to row 8.
$ws.Rows("8").Delete()
